$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-19 07:40:49"
$wsZhCn.Range("G4").Value = "2016-01-19 07:41:35"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-19 07:41:00"
$wsDeDe.Range("G4").Value = "2016-01-19 07:41:53"
